# Fruta / hortaliza, semanal
# Insert a new daily record as row 53 on the "Vega Monumental Concepción - Uva"
# sheet, pushing the existing rows 53-98 down to 54-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 53 (shifts rows 53:98 -> 54:99)
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new record's data
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = 'Vega Monumental Concepción'
$ws.Range("C53").Value = 'Bíobío'
$ws.Range("D53").Value = 44574
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 'Fruta'
$ws.Range("G53").Value = 100109
$ws.Range("H53").Value = 'Uva'
$ws.Range("I53").Value = 100109001
$ws.Range("J53").Value = 'Uva'
$ws.Range("K53").Value = 'Superior Seedless'
$ws.Range("L53").Value = 'Primera'
$ws.Range("M53").Value = 200
$ws.Range("N53").Value = 15000
$ws.Range("O53").Value = 16000
$ws.Range("P53").Value = 15500
$ws.Range("Q53").Value = '$/caja 15 kilos'
$ws.Range("R53").Value = 'Provincia de Limarí'
$ws.Range("S53").Value = 1033
$ws.Range("T53").Value = 15
